$wb = $excel.ActiveWorkbook

# The "optimization_parameters" sheet had a stray leftover row (row 16:
# "Sheet" / 3 / 4) that doesn't belong with the rest of the optimization
# parameters. Select it and delete the entire row, shifting everything
# below it up by one.
$ws = $wb.Worksheets.Item("optimization_parameters")
$ws.Rows.Item(16).Select()
$ws.Rows.Item(16).Delete()

# Finish up on the last sheet (optimization_diagnostics), which becomes
# the active/selected sheet in the saved workbook.
$last = $wb.Worksheets.Item("optimization_diagnostics")
$last.Activate()
